$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 62500680
$ws.Range("I9").Value = 524.5833
$ws.Range("J9").Value = 250001140
$ws.Range("K9").Value = 524.5833
$ws.Range("L9").Value = 250001140
$ws.Range("M9").Value = -355.5833
$ws.Range("N9").Value = -250001478

$ws.Range("H15").Value = 314324.78
$ws.Range("I15").Value = 314324.78
$ws.Range("K15").Value = 942974.3400000001
$ws.Range("M15").Value = -942805.3400000001

$ws.Range("H40").Value = 6830.294
$ws.Range("I40").Value = 4559.8184
$ws.Range("K40").Value = 4559.8184
$ws.Range("M40").Value = -4384.8184

$ws.Range("H58").Value = 803.5714
$ws.Range("I58").Value = 565
$ws.Range("K58").Value = 1695
$ws.Range("M58").Value = -1545

$ws.Range("H62").Value = 12409.875
$ws.Range("I62").Value = 12409.875
$ws.Range("K62").Value = 12409.875
$ws.Range("M62").Value = -11785.875

$ws.Range("H65").Value = 12409.875
$ws.Range("I65").Value = 12409.875
$ws.Range("K65").Value = 62049.375
$ws.Range("M65").Value = -58929.375

$ws.Range("H86").Value = 4638.5
$ws.Range("I86").Value = 3531.7646
$ws.Range("J86").Value = 6085.769
$ws.Range("K86").Value = 3531.7646
$ws.Range("L86").Value = 6085.769
$ws.Range("M86").Value = -2408.7646
$ws.Range("N86").Value = -8331.769

$ws.Range("H89").Value = 4638.5
$ws.Range("I89").Value = 3531.7646
$ws.Range("J89").Value = 6085.769
$ws.Range("K89").Value = 17658.823
$ws.Range("L89").Value = 30428.845
$ws.Range("M89").Value = -12042.823
$ws.Range("N89").Value = -41660.845

$ws.Range("H106").Value = 4275.5
$ws.Range("I106").Value = 3725.1538
$ws.Range("K106").Value = 3725.1538
$ws.Range("M106").Value = -3094.1538

$ws.Range("H112").Value = 4823.2354
$ws.Range("J112").Value = 4823.2354
$ws.Range("L112").Value = 14469.7062
$ws.Range("N112").Value = -16685.7062

$ws.Range("H127").Value = 2332.1904
$ws.Range("I127").Value = 836.8570999999999
$ws.Range("J127").Value = 5322.857
$ws.Range("K127").Value = 2510.5713
$ws.Range("L127").Value = 15968.571
$ws.Range("M127").Value = 2449.4287
$ws.Range("N127").Value = -25888.571

$ws.Range("H130").Value = 46500
$ws.Range("J130").Value = 46500
$ws.Range("L130").Value = 46500
$ws.Range("N130").Value = -56540

$ws.Range("H137").Value = 4128.778
$ws.Range("I137").Value = 3671.5715
$ws.Range("J137").Value = 5729
$ws.Range("K137").Value = 11014.7145
$ws.Range("L137").Value = 17187
$ws.Range("M137").Value = -8464.7145
$ws.Range("N137").Value = -22287

$ws.Range("H138").Value = 7181.183
$ws.Range("I138").Value = 7364.1665
$ws.Range("J138").Value = 7149.8145
$ws.Range("K138").Value = 22092.4995
$ws.Range("L138").Value = 21449.4435
$ws.Range("M138").Value = -16952.4995
$ws.Range("N138").Value = -31729.4435

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2458.1082
$ws.Range("I2").Value = 1967.2413
$ws.Range("J2").Value = 4237.5
$ws.Range("K2").Value = 1967.2413
$ws.Range("L2").Value = 4237.5
$ws.Range("M2").Value = -1854.2413
$ws.Range("N2").Value = -4463.5

$ws.Range("H32").Value = 18911.346
$ws.Range("I32").Value = 17631.96
$ws.Range("K32").Value = 17631.96
$ws.Range("M32").Value = -17344.96

$ws.Range("H44").Value = 69657.836
$ws.Range("J44").Value = 69657.836
$ws.Range("L44").Value = 69657.836
$ws.Range("N44").Value = -70633.836

$ws.Range("H45").Value = 4183.421
$ws.Range("I45").Value = 3285.4443
$ws.Range("J45").Value = 4991.6
$ws.Range("K45").Value = 3285.4443
$ws.Range("L45").Value = 4991.6
$ws.Range("M45").Value = -2908.4443
$ws.Range("N45").Value = -5745.6

$ws.Range("H61").Value = 18932.428
$ws.Range("I61").Value = 18902.6
$ws.Range("K61").Value = 18902.6
$ws.Range("M61").Value = -18690.6

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H74").Value = 6075.5884
$ws.Range("I74").Value = 4867.769
$ws.Range("K74").Value = 4867.769
$ws.Range("M74").Value = -3993.769

$ws.Range("H77").Value = 6075.5884
$ws.Range("I77").Value = 4867.769
$ws.Range("K77").Value = 24338.845
$ws.Range("M77").Value = -19970.845

$ws.Range("H97").Value = 1585.1
$ws.Range("I97").Value = 1658
$ws.Range("J97").Value = 200
$ws.Range("K97").Value = 1658
$ws.Range("L97").Value = 200
$ws.Range("M97").Value = -1162
$ws.Range("N97").Value = -1192

$ws.Range("H102").Value = 10000
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws.Range("H110").Value = 1697.3334
$ws.Range("I110").Value = 1453.6666
$ws.Range("K110").Value = 1453.6666
$ws.Range("M110").Value = 591.3334

$ws.Range("H116").Value = 2458.1082
$ws.Range("I116").Value = 1967.2413
$ws.Range("J116").Value = 4237.5
$ws.Range("K116").Value = 1967.2413
$ws.Range("L116").Value = 4237.5
$ws.Range("M116").Value = 326.7587000000001
$ws.Range("N116").Value = -8825.5

$ws.Range("H122").Value = 8416.959999999999
$ws.Range("I122").Value = 7842.6665
$ws.Range("J122").Value = 11432
$ws.Range("K122").Value = 23527.9995
$ws.Range("L122").Value = 34296
$ws.Range("M122").Value = -21077.9995
$ws.Range("N122").Value = -39196

$ws.Range("H136").Value = 18932.428
$ws.Range("I136").Value = 18902.6
$ws.Range("K136").Value = 56707.8
$ws.Range("M136").Value = -54157.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2458.1082
$ws.Range("I3").Value = 1967.2413
$ws.Range("J3").Value = 4237.5
$ws.Range("K3").Value = 1967.2413
$ws.Range("L3").Value = 4237.5
$ws.Range("M3").Value = -1853.2413
$ws.Range("N3").Value = -4465.5

$ws.Range("H20").Value = 2577.889
$ws.Range("I20").Value = 2376.647
$ws.Range("K20").Value = 2376.647
$ws.Range("M20").Value = -2129.647

$ws.Range("H35").Value = 116000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 116000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 116000
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -116620

$ws.Range("H86").Value = 9998.4375
$ws.Range("I86").Value = 10887.111
$ws.Range("J86").Value = 8855.857
$ws.Range("K86").Value = 10887.111
$ws.Range("L86").Value = 8855.857
$ws.Range("M86").Value = -9764.111000000001
$ws.Range("N86").Value = -11101.857

$ws.Range("H89").Value = 9998.4375
$ws.Range("I89").Value = 10887.111
$ws.Range("J89").Value = 8855.857
$ws.Range("K89").Value = 54435.55500000001
$ws.Range("L89").Value = 44279.285
$ws.Range("M89").Value = -48819.55500000001
$ws.Range("N89").Value = -55511.285

$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").ClearContents()

$ws.Range("H105").Value = 3698.3914
$ws.Range("I105").Value = 3769
$ws.Range("J105").Value = 3498.3333
$ws.Range("K105").Value = 3769
$ws.Range("L105").Value = 3498.3333
$ws.Range("M105").Value = -2022
$ws.Range("N105").Value = -6992.3333

$ws.Range("H107").Value = 3616
$ws.Range("I107").Value = 3616
$ws.Range("K107").Value = 3616
$ws.Range("M107").Value = -1696

$ws.Range("H134").Value = 6306.6665
$ws.Range("I134").Value = 5333.483
$ws.Range("K134").Value = 16000.449
$ws.Range("M134").Value = -13465.449

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 2785.8538
$ws.Range("J7").Value = 2716.1875
$ws.Range("L7").Value = 2716.1875
$ws.Range("N7").Value = -2942.1875

$ws.Range("H22").Value = 2989.6
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 4582.6665
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 4582.6665
$ws.Range("M22").Value = -250
$ws.Range("N22").Value = -5282.6665

$ws.Range("H99").Value = 5413.355
$ws.Range("I99").Value = 6221.875
$ws.Range("J99").Value = 2641.2856
$ws.Range("K99").Value = 6221.875
$ws.Range("L99").Value = 2641.2856
$ws.Range("M99").Value = -4723.875
$ws.Range("N99").Value = -5637.2856

$ws.Range("H107").Value = 2228.8572
$ws.Range("I107").Value = 2383.8823
$ws.Range("J107").Value = 1989.2727
$ws.Range("K107").Value = 2383.8823
$ws.Range("L107").Value = 1989.2727
$ws.Range("M107").Value = -463.8823000000002
$ws.Range("N107").Value = -5829.2727

$ws.Range("H122").Value = 7999.75
$ws.Range("I122").Value = 7000
$ws.Range("J122").Value = 8999.5
$ws.Range("K122").Value = 21000
$ws.Range("L122").Value = 26998.5
$ws.Range("M122").Value = -18550
$ws.Range("N122").Value = -31898.5

$ws.Range("H126").Value = 5413.355
$ws.Range("I126").Value = 6221.875
$ws.Range("J126").Value = 2641.2856
$ws.Range("K126").Value = 18665.625
$ws.Range("L126").Value = 7923.8568
$ws.Range("M126").Value = -16195.625
$ws.Range("N126").Value = -12863.8568

$ws.Range("H141").Value = 311577.38
$ws.Range("I141").Value = 89832.664
$ws.Range("J141").Value = 394731.62
$ws.Range("K141").Value = 89832.664
$ws.Range("L141").Value = 394731.62
$ws.Range("M141").Value = -84652.664
$ws.Range("N141").Value = -405091.62

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 189.875
$ws.Range("I2").Value = 239.83333
$ws.Range("K2").Value = 1438.99998
$ws.Range("M2").Value = -1325.99998

$ws.Range("H21").Value = 35
$ws.Range("I21").Value = 50
$ws.Range("J21").Value = 5
$ws.Range("K21").Value = 150
$ws.Range("L21").Value = 15
$ws.Range("M21").Value = 23
$ws.Range("N21").Value = -361

$ws.Range("H107").Value = 555.4
$ws.Range("J107").Value = 628.25
$ws.Range("L107").Value = 1884.75
$ws.Range("N107").Value = -5724.75

$ws.Range("H109").Value = 3668.8
$ws.Range("I109").Value = 3668.8
$ws.Range("K109").Value = 11006.4
$ws.Range("M109").Value = -9966.400000000001

$ws.Range("H113").Value = 1894.875
$ws.Range("J113").Value = 1895.4
$ws.Range("L113").Value = 5686.200000000001
$ws.Range("N113").Value = -10026.2

$ws.Range("H121").Value = 703946.4399999999
$ws.Range("J121").Value = 1266861.8
$ws.Range("L121").Value = 3800585.4
$ws.Range("N121").Value = -3803205.4

$ws.Range("H131").Value = 2886.5
$ws.Range("J131").Value = 4000
$ws.Range("L131").Value = 12000
$ws.Range("N131").Value = -22080

$ws.Range("H136").Value = 3299.5
$ws.Range("I136").Value = 3299.5
$ws.Range("K136").Value = 9898.5
$ws.Range("M136").Value = -4798.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6076.4443
$ws.Range("I70").Value = 4849
$ws.Range("J70").Value = 6427.143
$ws.Range("K70").Value = 4849
$ws.Range("L70").Value = 6427.143
$ws.Range("M70").Value = -4579
$ws.Range("N70").Value = -6967.143

$ws.Range("H73").Value = 6076.4443
$ws.Range("I73").Value = 4849
$ws.Range("J73").Value = 6427.143
$ws.Range("K73").Value = 4849
$ws.Range("L73").Value = 6427.143
$ws.Range("M73").Value = -3913
$ws.Range("N73").Value = -8299.143

$ws.Range("H80").Value = 4000
$ws.Range("I80").Value = 4000
$ws.Range("K80").Value = 4000
$ws.Range("M80").Value = -3002

$ws.Range("H83").Value = 4000
$ws.Range("I83").Value = 4000
$ws.Range("K83").Value = 20000
$ws.Range("M83").Value = -15008

$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()

$ws.Range("H102").Value = 3527.1667
$ws.Range("I102").Value = 3112.1035
$ws.Range("K102").Value = 3112.1035
$ws.Range("M102").Value = -1490.1035

$ws.Range("H126").Value = 11813.909
$ws.Range("I126").Value = 11229.941
$ws.Range("K126").Value = 33689.823
$ws.Range("M126").Value = -31219.823

$ws.Range("H132").Value = 8632.1875
$ws.Range("I132").Value = 7909.05
$ws.Range("J132").Value = 12247.875
$ws.Range("K132").Value = 23727.15
$ws.Range("L132").Value = 36743.625
$ws.Range("M132").Value = -21197.15
$ws.Range("N132").Value = -41803.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3668.8
$ws.Range("I7").Value = 3446.3333
$ws.Range("K7").Value = 3446.3333
$ws.Range("M7").Value = -3334.3333

$ws.Range("H40").Value = 5821.7856
$ws.Range("I40").Value = 4885
$ws.Range("J40").Value = 18000
$ws.Range("K40").Value = 4885
$ws.Range("L40").Value = 18000
$ws.Range("M40").Value = -4749
$ws.Range("N40").Value = -18272

$ws.Range("H43").Value = 39999
$ws.Range("I43").Value = 39999
$ws.Range("K43").Value = 39999
$ws.Range("M43").Value = -39806

$ws.Range("H46").Value = 4401.619
$ws.Range("J46").Value = 4289.3
$ws.Range("L46").Value = 4289.3
$ws.Range("N46").Value = -4665.3

$ws.Range("H55").Value = 857.55554
$ws.Range("J55").Value = 1998
$ws.Range("L55").Value = 1998
$ws.Range("N55").Value = -2344

$ws.Range("H64").Value = 215000
$ws.Range("J64").Value = 215000
$ws.Range("L64").Value = 215000
$ws.Range("N64").Value = -215450

$ws.Range("H67").Value = 215000
$ws.Range("J67").Value = 215000
$ws.Range("L67").Value = 215000
$ws.Range("N67").Value = -216560

$ws.Range("H68").Value = 10907
$ws.Range("I68").Value = 7998.8335
$ws.Range("K68").Value = 7998.8335
$ws.Range("M68").Value = -7249.8335

$ws.Range("H70").Value = 29999
$ws.Range("J70").Value = 29999
$ws.Range("L70").Value = 29999
$ws.Range("N70").Value = -30539

$ws.Range("H71").Value = 10907
$ws.Range("I71").Value = 7998.8335
$ws.Range("K71").Value = 39994.1675
$ws.Range("M71").Value = -36250.1675

$ws.Range("H73").Value = 29999
$ws.Range("J73").Value = 29999
$ws.Range("L73").Value = 29999
$ws.Range("N73").Value = -31871

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H82").Value = 1172.5
$ws.Range("I82").Value = 896.6667
$ws.Range("K82").Value = 896.6667
$ws.Range("M82").Value = -535.6667

$ws.Range("H85").Value = 1172.5
$ws.Range("I85").Value = 896.6667
$ws.Range("K85").Value = 896.6667
$ws.Range("M85").Value = 351.3333

$ws.Range("H100").Value = 6433.7666
$ws.Range("J100").Value = 6596.2144
$ws.Range("L100").Value = 6596.2144
$ws.Range("N100").Value = -7678.2144

$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H122").Value = 4832.8335
$ws.Range("I122").Value = 3920.4211
$ws.Range("J122").Value = 8300
$ws.Range("K122").Value = 11761.2633
$ws.Range("L122").Value = 24900
$ws.Range("M122").Value = -9311.263300000001
$ws.Range("N122").Value = -29800

$ws.Range("H126").Value = 3668.8
$ws.Range("I126").Value = 3446.3333
$ws.Range("K126").Value = 10338.9999
$ws.Range("M126").Value = -7868.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 12232.417
$ws.Range("J45").Value = 11051.223
$ws.Range("L45").Value = 11051.223
$ws.Range("N45").Value = -12033.223

$ws.Range("H62").Value = 8000
$ws.Range("J62").Value = 8000
$ws.Range("L62").Value = 8000
$ws.Range("N62").Value = -9248

$ws.Range("H65").Value = 8000
$ws.Range("J65").Value = 8000
$ws.Range("L65").Value = 40000
$ws.Range("N65").Value = -46240

$ws.Range("H122").Value = 5789.4116
$ws.Range("I122").Value = 5624.9165
$ws.Range("K122").Value = 16874.7495
$ws.Range("M122").Value = -14424.7495

$ws.Range("H126").Value = 1968.3405
$ws.Range("I126").Value = 1570.25
$ws.Range("J126").Value = 4243.143
$ws.Range("K126").Value = 4710.75
$ws.Range("L126").Value = 12729.429
$ws.Range("M126").Value = -2240.75
$ws.Range("N126").Value = -17669.429

$ws.Range("H136").Value = 3270.2104
$ws.Range("I136").Value = 2455
$ws.Range("J136").Value = 5897
$ws.Range("K136").Value = 7365
$ws.Range("L136").Value = 17691
$ws.Range("M136").Value = -4815
$ws.Range("N136").Value = -22791
